$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.557.68'
$ws.Range('E2').Value = '  +3.08%  '
$ws.Range('D3').Value = '1.598.25'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '26.79'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.40'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').Value = '1.826.08'
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '1.607.82'
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('D15').Value = '29.561.80'
$ws.Range('E15').Value = '  +3.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.536'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').Value = '0.0₃0693'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.09%  '
$ws.Range('E28').Value = '  +1.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.78%  '
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0476'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.51%  '
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('D34').Value = '1.435.09'
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.85%  '
$ws.Range('E38').Value = '  +3.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.538'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0493'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.73%  '
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '53.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +26.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.800'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +20.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.30%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '1.736.90'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.87%  '
